$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new log row (row 9) for the mail processed at 2025-08-18 21:12:41
$ws.Range("A9").Value = "Geen onderwerp"
$ws.Range("B9").Value = "onbekend"
$ws.Range("D9").Value = "Overig"
$ws.Range("F9").Value = "2025-08-18 21:12:41"
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = $true
$ws.Range("I9").Value = $false
$ws.Range("J9").Value = $false

# Extend the conditional formatting ranges so they keep covering the data
# (they previously stopped at row 8, now the table goes to row 9)
$colLetters = @("D", "G", "H", "I", "J")
foreach ($col in $colLetters) {
    $oldRange = $ws.Range($col + "2:" + $col + "8")
    $newRange = $ws.Range($col + "2:" + $col + "9")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary: "Overig" category count goes from 1 to 2
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 2
